# Update SwaadSutra_Daily_2026-01-20.xlsx
# - Insert a new order row at the top of "Daily Orders" (row 2), pushing
#   existing orders down by one row (order #24, Udita Roy @ 15:03).
# - Refresh the "Summary" totals.
# - Refresh the "Items Breakdown" counts for "Appe Chutney".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Daily Orders sheet: insert a new row for the newest order (#24)
# ---------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Daily Orders")

# Shift existing data rows (2..5) down to (3..6) before filling the new row.
$orders.Rows.Item(2).Insert()

$newRow = 2
$orders.Cells.Item($newRow, 1).Value = 24
$orders.Cells.Item($newRow, 2).Value = "2026-01-20 15:03"
$orders.Cells.Item($newRow, 3).Value = "Udita Roy"
$orders.Cells.Item($newRow, 4).Value = "A-1603"

# Phone numbers must stay text (not be coerced into a number) - format the
# cell as Text first, assign, then drop the formatting so no style sticks.
$phoneCell = $orders.Cells.Item($newRow, 5)
$phoneCell.NumberFormat = "@"
$phoneCell.Value = "7061856166"
$phoneCell.ClearFormats()

$orders.Cells.Item($newRow, 6).Value = "Appe Chutney x1"
$orders.Cells.Item($newRow, 7).Value = 60
$orders.Cells.Item($newRow, 8).Value = "NEW"
$orders.Cells.Item($newRow, 9).Value = "PENDING"

# Collection date must stay as plain text "2026-01-21", not an Excel date.
$dateCell = $orders.Cells.Item($newRow, 10)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-01-21"
$dateCell.ClearFormats()

$orders.Cells.Item($newRow, 11).Value = "09:30"
$orders.Cells.Item($newRow, 12).Value = "Less spicy. Flavourful"
$orders.Cells.Item($newRow, 13).Value = ""
$orders.Cells.Item($newRow, 14).Value = ""

# ---------------------------------------------------------------------
# 2) Summary sheet: bump Total Orders, New count, and Total Revenue
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(2, 1).Value = 5    # Total Orders: 4 -> 5
$summary.Cells.Item(2, 2).Value = 2    # New: 1 -> 2
$summary.Cells.Item(2, 7).Value = 320  # Total Revenue: 260 -> 320

# ---------------------------------------------------------------------
# 3) Items Breakdown sheet: bump Appe Chutney quantity and revenue
# ---------------------------------------------------------------------
$items = $wb.Worksheets.Item("Items Breakdown")
$items.Cells.Item(2, 2).Value = 3      # Quantity Ordered: 2 -> 3
$items.Cells.Item(2, 3).Value = 180    # Revenue: 120 -> 180
